$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New config rows (MicrosoftEdge / ZohoConnect) added under the Settings sheet
$ws.Range("A6").Value = "MicrosoftEdge"
$ws.Range("B6").Value = "msedge"
$ws.Range("C6").Value = "Close the application"

$ws.Range("B7").Value = "https://www.zoho.com/connect/"
$ws.Range("A7").Value = "ZohoConnect"
$ws.Range("C7").Value = "Close Tab"

# Hyperlink the Zoho Connect URL cell
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.zoho.com/connect/") | Out-Null

# Settings becomes the active/selected sheet and tab
$ws.Activate()
$ws.Range("C9").Select() | Out-Null
